# "Added Model Convergence Test"
#
# case12.xlsx : Bus / Gen / Branch sheets for a power-flow model.
#  - Bus!J (baseKV) is rebased from a 12.4kV system to a 240V system
#    (column header + all data values), and a new comment documenting
#    the new sBase/vBase/iBase/zBase is dropped on Bus!J1.
#  - Branch!E (x, per-unit reactance) is recomputed for the new base
#    (zBase = 0.576) and the existing Branch!E1 comment gets the
#    "Sweetbriar" resistance derivation appended.
#  - Final active sheet/selection moves from Gen to Branch.

$wb = $excel.ActiveWorkbook

$busSheet    = $wb.Worksheets.Item("Bus")
$genSheet    = $wb.Worksheets.Item("Gen")
$branchSheet = $wb.Worksheets.Item("Branch")

# ---------------------------------------------------------------------
# Bus sheet: baseKV -> baseV, values 12.4 -> 240, new header comment
# ---------------------------------------------------------------------
$busSheet.Activate()

$busSheet.Range("J1").Value = "baseV"

$busSheet.Range("J2").Value = 240
$busSheet.Range("J3").Value = 240
$busSheet.Range("J4").Value = 240
$busSheet.Range("J5").Value = 240
$busSheet.Range("J6").Value = 240
$busSheet.Range("J7").Value = 240
$busSheet.Range("J8").Value = 240
$busSheet.Range("J9").Value = 240
$busSheet.Range("J10").Value = 240

$baseNote = "sBase 100 kVA`n 100000 VA`nvBase 240 V`niBase 416.6667 `nzBase 0.576 `n"
$busComment = $busSheet.Range("J1").AddComment($baseNote)

$busSheet.Range("H17").Select()

# ---------------------------------------------------------------------
# Gen sheet: no data changes, selection only moves to G1
# ---------------------------------------------------------------------
$genSheet.Activate()
$genSheet.Range("G1").Select()

# ---------------------------------------------------------------------
# Branch sheet: recompute per-unit reactance, extend comment, page setup
# ---------------------------------------------------------------------
$branchSheet.Activate()

$branchSheet.Range("E2").Value = 0.02192
$branchSheet.Range("E3").Value = 0.02192
$branchSheet.Range("E4").Value = 0.02192
$branchSheet.Range("E5").Value = 0.02192
$branchSheet.Range("E6").Value = 0.02192
$branchSheet.Range("E7").Value = 0.02192
$branchSheet.Range("E8").Value = 0.02192
$branchSheet.Range("E9").Value = 0.02192

$reactanceComment = $branchSheet.Range("E1").Comment
$reactanceText = "Per-Unit Reactance`nSweetbriar has resistance of 0.101 / 1000ft. Assume 125ft per cable = 0.012625`nzBase = 0.576`npu_X = 0.012625/zBase = 0.02192"
$reactanceComment.Text($reactanceText)

# Bold the heading / the " " separator run, matching the source formatting
# (best effort - engine may not persist per-run formatting on comments).
$reactanceFrame = $reactanceComment.Shape.TextFrame
$reactanceFrame.Characters(1, 19).Font.Bold = $true
$reactanceFrame.Characters(50, 1).Font.Bold = $true

$branchSheet.PageSetup.Orientation = 1

$branchSheet.Range("E2:E9").Select()

# ---------------------------------------------------------------------
# Final selected tab is Branch (activeTab = 2 / 0-based 3rd sheet)
# ---------------------------------------------------------------------
$branchSheet.Activate()
